# Generate Report for Handback
# - Mark zh-cn / de-de rows as "Handed back: in sync with en-US" (was "In Translation")
# - de-de handback completed: link target .md file + target .xlf file + handback datetime
# - zh-cn not handed back yet: just refresh the "no handback" placeholder datetime
# - Widen columns that now hold the longer status text / file-name links

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: status cells for both locales, both rows
# ---------------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen the zh-cn / de-de status columns on the Overview sheet
$overview.Columns.Item(5).ColumnWidth = 29.1443
$overview.Columns.Item(6).ColumnWidth = 29.1443

# ---------------------------------------------------------------------------
# zh-cn sheet: status text, target/handback-file links not filled in yet
# (no successful handback yet), just refresh the "never handed back" datetime
# ---------------------------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$zhcn.Range("I2").Value = "58c5f947-0739-4c0a-96ab-02c3696caeda.md"
$zhcn.Range("J2").Value = "58c5f947-0739-4c0a-96ab-02c3696caeda.e011d484fe5af3c2137178b777831a1c934b373c.zh-cn.xlf"
$zhcn.Range("I3").Value = "69ce8d27-da3f-45e6-8bd0-27248c73a872.md"
$zhcn.Range("J3").Value = "69ce8d27-da3f-45e6-8bd0-27248c73a872.f1b2c7ed9ef4f1f980cac7acd748d532a6e515d3.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-08-15 12:20:47"
$zhcn.Range("K3").Value = "2016-08-15 12:20:47"

$zhcn.Columns.Item(3).ColumnWidth = 29.1443
$zhcn.Columns.Item(9).ColumnWidth = 39.1667
$zhcn.Columns.Item(10).ColumnWidth = 39.1667

# rebuild hyperlinks in document order: A2, I2 (new), A3, I3 (new)
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20dbf1d5f6669380c6adb40f0d806eac3c9c0725/e2e/58c5f947-0739-4c0a-96ab-02c3696caeda.md", "", "", "58c5f947-0739-4c0a-96ab-02c3696caeda.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20dbf1d5f6669380c6adb40f0d806eac3c9c0725/e2e/58c5f947-0739-4c0a-96ab-02c3696caeda.md", "", "", "58c5f947-0739-4c0a-96ab-02c3696caeda.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20dbf1d5f6669380c6adb40f0d806eac3c9c0725/e2e/69ce8d27-da3f-45e6-8bd0-27248c73a872.md", "", "", "69ce8d27-da3f-45e6-8bd0-27248c73a872.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20dbf1d5f6669380c6adb40f0d806eac3c9c0725/e2e/69ce8d27-da3f-45e6-8bd0-27248c73a872.md", "", "", "69ce8d27-da3f-45e6-8bd0-27248c73a872.md")

foreach ($addr in @("A2", "I2", "A3", "I3")) {
    $r = $zhcn.Range($addr)
    $r.Font.Underline = $true
    $r.Font.Color = 15570276
}

# ---------------------------------------------------------------------------
# de-de sheet: status text, AND handback fully completed - fill target file,
# handback file, handback datetime
# ---------------------------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("I2").Value = "58c5f947-0739-4c0a-96ab-02c3696caeda.md"
$dede.Range("J2").Value = "58c5f947-0739-4c0a-96ab-02c3696caeda.e011d484fe5af3c2137178b777831a1c934b373c.de-de.xlf"
$dede.Range("I3").Value = "69ce8d27-da3f-45e6-8bd0-27248c73a872.md"
$dede.Range("J3").Value = "69ce8d27-da3f-45e6-8bd0-27248c73a872.f1b2c7ed9ef4f1f980cac7acd748d532a6e515d3.de-de.xlf"

$dede.Range("K2").Value = "2016-08-15 12:20:54"
$dede.Range("K3").Value = "2016-08-15 12:20:54"

$dede.Columns.Item(3).ColumnWidth = 29.1443
$dede.Columns.Item(9).ColumnWidth = 39.1667
$dede.Columns.Item(10).ColumnWidth = 39.1667

# rebuild hyperlinks in document order: A2, I2 (new), A3, I3 (new)
$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20dbf1d5f6669380c6adb40f0d806eac3c9c0725/e2e/58c5f947-0739-4c0a-96ab-02c3696caeda.md", "", "", "58c5f947-0739-4c0a-96ab-02c3696caeda.md")
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20dbf1d5f6669380c6adb40f0d806eac3c9c0725/e2e/58c5f947-0739-4c0a-96ab-02c3696caeda.md", "", "", "58c5f947-0739-4c0a-96ab-02c3696caeda.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20dbf1d5f6669380c6adb40f0d806eac3c9c0725/e2e/69ce8d27-da3f-45e6-8bd0-27248c73a872.md", "", "", "69ce8d27-da3f-45e6-8bd0-27248c73a872.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/20dbf1d5f6669380c6adb40f0d806eac3c9c0725/e2e/69ce8d27-da3f-45e6-8bd0-27248c73a872.md", "", "", "69ce8d27-da3f-45e6-8bd0-27248c73a872.md")

foreach ($addr in @("A2", "I2", "A3", "I3")) {
    $r = $dede.Range($addr)
    $r.Font.Underline = $true
    $r.Font.Color = 15570276
}

Write-Output "Handback report generated"
